# Angular - 01.INTRO TO ANGULAR AND TYPESCRIPT - Complete
# Inserts a new row (7) describing compiling to ES5 with tsc, and a new
# column (D) with an additional tsc invocation example on row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 7 - this shifts the former rows 7-15
#    down to 8-16 and keeps their formatting (style, merges) intact.
$ws.Rows.Item(7).Insert()

# 2. New cell D6 - extra tsc invocation example next to the existing
#    "To compile your code" / "tsc myfile.ts" / "Compilation output ..." row.
#    Apply the shared row style first so the cell keeps s="1" like its
#    A6:C6 neighbours, then set its value.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D6").Value = ".\node_modules\.bin\tsc index.ts"

# 3. Populate the newly inserted row 7 (order matters so shared strings end
#    up in the same slots as the source workbook).
$ws.Range("B7").Value = " .\node_modules\.bin\tsc -t es5 6.KeyValuePairs.ts"
$ws.Range("A7").Value = "To compile your code with for target version of es 5"
$ws.Range("C7").Value = "Compilation output is plain JavaScript"

# Row 7 mirrors the row-height of the other wrapped, 3-line rows (45pt).
$ws.Rows.Item(7).RowHeight = 45

# D7 stays empty but still carries the same cell style as the rest of the
# row (style index 1) - copy formats only from a cell that already uses it.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 4. The row that used to be r="15" (now r="16") had an explicit 60pt
#    height; after the insert it keeps that value, but the edited workbook
#    trims it back down to 30pt.
$ws.Rows.Item(16).RowHeight = 30

# 5. New column D needs an explicit width (~52.57 characters wide).
$ws.Columns.Item(4).ColumnWidth = 51.65

# 6. Match the author's final selection.
$ws.Range("C7").Select() | Out-Null
